$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.235.09'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '2.381.80'
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("E4").Value = '  +0.77%  '
$ws.Range("D5").Value = '562.53'
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").Value = '138.86'
$ws.Range("E6").Value = '  -0.44%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '0.531'
$ws.Range("E8").Value = '  +1.03%  '
$ws.Range("D9").Value = '2.381.47'
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("E11").Value = '  -1.18%  '
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").Value = '2.836.74'
$ws.Range("E15").Value = '  +1.32%  '
$ws.Range("D16").Value = '0.0000167'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '60.085.69'
$ws.Range("E17").Value = '  -0.92%  '
$ws.Range("D18").Value = '2.381.99'
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("D19").Value = '8.11'
$ws.Range("E19").Value = '  +11.26%  '
$ws.Range("D20").Value = '10.53'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = '322.83'
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '1.83'
$ws.Range("E25").Value = '  -2.58%  '
$ws.Range("D26").Value = '64.24'
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("D27").Value = '561.88'
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("D28").Value = '8.09'
$ws.Range("E28").Value = '  -5.70%  '
$ws.Range("D29").Value = '2.497.46'
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("D30").Value = '0.0₃0925'
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("E33").Value = '  -1.73%  '
$ws.Range("D34").Value = '0.132'
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("D36").Value = '1.45'
$ws.Range("E36").Value = '  +5.62%  '
$ws.Range("D37").Value = '154.20'
$ws.Range("E37").Value = '  +4.89%  '
$ws.Range("E38").Value = '  -0.60%  '
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("D40").Value = '18.18'
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").Value = '5.08'
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = '41.65'
$ws.Range("E43").Value = '  +1.70%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '1.67'
$ws.Range("E44").Value = '  +0.82%  '
$ws.Range("E45").Value = '  +4.27%  '
$ws.Range("E46").Value = '  +2.94%  '
$ws.Range("D47").Value = '140.42'
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").Value = '3.53'
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("E49").Value = '  +0.77%  '
$ws.Range("D50").Value = '0.0503'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = '19.19'
$ws.Range("E51").Value = '  -0.70%  '
